$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "43.905.06"
$ws.Cells.Item(2, 5).Value = "  -1.13%  "
$ws.Cells.Item(3, 4).Value = "2.191.05"
$ws.Cells.Item(3, 5).Value = "  -2.35%  "
Set-TextValue 4 4 "1.01"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 5).Value = "  -4.23%  "
Set-TextValue 6 4 "87.78"
$ws.Cells.Item(6, 5).Value = "  -5.70%  "
Set-TextValue 7 4 "0.563"
$ws.Cells.Item(7, 5).Value = "  -1.28%  "
$ws.Cells.Item(8, 5).Value = "  -0.08%  "
$ws.Cells.Item(9, 5).Value = "  -8.63%  "
Set-TextValue 10 4 "32.04"
$ws.Cells.Item(10, 5).Value = "  -7.20%  "
Set-TextValue 11 4 "0.0761"
$ws.Cells.Item(11, 5).Value = "  -5.98%  "
$ws.Cells.Item(12, 5).Value = "  -1.81%  "
Set-TextValue 13 4 "6.72"
$ws.Cells.Item(13, 5).Value = "  -5.95%  "
$ws.Cells.Item(14, 4).Value = "2.523.82"
$ws.Cells.Item(14, 5).Value = "  -2.41%  "
$ws.Cells.Item(15, 4).Value = "2.248.02"
$ws.Cells.Item(15, 5).Value = "  +0.43%  "
Set-TextValue 16 4 "12.87"
$ws.Cells.Item(16, 5).Value = "  -5.19%  "
Set-TextValue 17 4 "0.762"
$ws.Cells.Item(17, 5).Value = "  -9.00%  "
$ws.Cells.Item(18, 4).Value = "43.431.87"
$ws.Cells.Item(18, 5).Value = "  -1.45%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0875"
$ws.Cells.Item(19, 5).Value = "  -8.86%  "
Set-TextValue 20 4 "5.79"
$ws.Cells.Item(20, 5).Value = "  -9.12%  "
Set-TextValue 21 4 "10.65"
$ws.Cells.Item(21, 5).Value = "  -13.51%  "
Set-TextValue 22 4 "62.48"
$ws.Cells.Item(22, 5).Value = "  -4.84%  "
Set-TextValue 23 4 "228.68"
$ws.Cells.Item(23, 5).Value = "  -3.56%  "
$ws.Cells.Item(24, 5).Value = "  -6.26%  "
$ws.Cells.Item(25, 5).Value = "  -0.11%  "
$ws.Cells.Item(26, 5).Value = "  -8.81%  "
$ws.Cells.Item(27, 5).Value = "  +0.36%  "
$ws.Cells.Item(28, 2).Value = "Cosmos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 28 4 "9.13"
$ws.Cells.Item(28, 5).Value = "  -7.01%  "
$ws.Cells.Item(29, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 29 4 "35.10"
$ws.Cells.Item(29, 5).Value = "  -8.91%  "
Set-TextValue 30 4 "18.99"
$ws.Cells.Item(30, 5).Value = "  -5.27%  "
Set-TextValue 31 4 "145.23"
$ws.Cells.Item(31, 5).Value = "  -5.45%  "
Set-TextValue 32 4 "5.25"
$ws.Cells.Item(32, 5).Value = "  -11.35%  "
Set-TextValue 33 4 "2.51"
$ws.Cells.Item(33, 5).Value = "  -5.48%  "
$ws.Cells.Item(34, 5).Value = "  -9.56%  "
$ws.Cells.Item(35, 5).Value = "  -3.63%  "
Set-TextValue 36 4 "2.85"
$ws.Cells.Item(36, 5).Value = "  -8.22%  "
Set-TextValue 37 4 "0.101"
$ws.Cells.Item(37, 5).Value = "  -6.96%  "
Set-TextValue 38 4 "1.64"
$ws.Cells.Item(38, 5).Value = "  -9.94%  "
$ws.Cells.Item(39, 2).Value = "Celestia"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue 39 4 "13.13"
$ws.Cells.Item(39, 5).Value = "  -9.56%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 40 4 "0.0277"
$ws.Cells.Item(40, 5).Value = "  -7.60%  "
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 41 4 "3.47"
$ws.Cells.Item(41, 5).Value = "  -8.95%  "
$ws.Cells.Item(42, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 42 4 "1.00"
$ws.Cells.Item(42, 5).Value = "  -0.37%  "
Set-TextValue 43 4 "3.03"
$ws.Cells.Item(43, 5).Value = "  -12.04%  "
$ws.Cells.Item(44, 4).Value = "1.741.99"
$ws.Cells.Item(44, 5).Value = "  +0.38%  "
Set-TextValue 45 4 "1.61"
$ws.Cells.Item(45, 5).Value = "  +1.50%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 46 4 "14.09"
$ws.Cells.Item(46, 5).Value = "  -1.00%  "
$ws.Cells.Item(47, 2).Value = "BitcoinSV"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue 47 4 "72.35"
$ws.Cells.Item(47, 5).Value = "  -9.92%  "
$ws.Cells.Item(48, 2).Value = "ordi"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue 48 4 "65.71"
$ws.Cells.Item(48, 5).Value = "  -5.16%  "
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 49 4 "0.170"
$ws.Cells.Item(49, 5).Value = "  -11.67%  "
$ws.Cells.Item(50, 2).Value = "HuobiToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 50 4 "2.73"
$ws.Cells.Item(50, 5).Value = "  +8.39%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 51 4 "90.73"
$ws.Cells.Item(51, 5).Value = "  -8.79%  "
